$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text formatting on the numeric-looking columns (score, reviews count)
# before writing, then reset the style afterwards so no stray number format
# is left on the cells (matches the plain/default-styled source cells).
$scoreRange = $ws.Range("C2:C26")
$reviewsRange = $ws.Range("E2:E26")
$scoreRange.NumberFormat = "@"
$reviewsRange.NumberFormat = "@"

$ws.Range("A2").Value = 'Cosy place in Paris'
$ws.Range("B2").Value = 'US$2,396'
$ws.Range("C2").Value = '8.0'
$ws.Range("D2").Value = 'Very Good'
$ws.Range("E2").Value = '1'

$ws.Range("A3").Value = '26 Faubourg - Ex-Hotel de Reims'
$ws.Range("B3").Value = 'US$3,327'
$ws.Range("C3").Value = '7.9'
$ws.Range("D3").Value = 'Good'
$ws.Range("E3").Value = '822'

$ws.Range("A4").Value = 'Luxurious Flat next to Eiffel Tower (WiFi/Netflix)'
$ws.Range("B4").Value = 'US$23,076'
$ws.Range("C4").Value = '8.8'
$ws.Range("D4").Value = 'Excellent'
$ws.Range("E4").Value = '12'

$ws.Range("A5").Value = 'Hotel Migny Opéra Montmartre'
$ws.Range("B5").Value = 'US$6,005'
$ws.Range("C5").Value = '8.0'
$ws.Range("D5").Value = 'Very Good'
$ws.Range("E5").Value = '4,139'

$ws.Range("A6").Value = 'Hôtel Saint-Pétersbourg Opéra & Spa'
$ws.Range("B6").Value = 'US$20,034'
$ws.Range("C6").Value = '8.6'
$ws.Range("D6").Value = 'Excellent'
$ws.Range("E6").Value = '4,347'

$ws.Range("A7").Value = 'CMG Charonne / Voltaire III'
$ws.Range("B7").Value = 'US$9,289'
$ws.Range("C7").Value = '6.9'
$ws.Range("D7").Value = 'Review score'
$ws.Range("E7").Value = '78'

$ws.Range("A8").Value = 'Luxe Apartment 165m2 8pers Victor Hugo trocadero foch Champs Elysées'
$ws.Range("B8").Value = 'US$29,966'
$ws.Range("C8").Value = '7.3'
$ws.Range("D8").Value = 'Good'
$ws.Range("E8").Value = '9'

$ws.Range("A9").Value = 'Le berlier 13'
$ws.Range("B9").Value = 'US$8,582'
$ws.Range("C9").Value = '9.7'
$ws.Range("D9").Value = 'Exceptional'
$ws.Range("E9").Value = '7'

$ws.Range("A10").Value = 'FINESTATE Coliving Champs-Elysées'
$ws.Range("B10").Value = 'US$3,301'
$ws.Range("C10").Value = '9.1'
$ws.Range("D10").Value = 'Wonderful'
$ws.Range("E10").Value = '25'

$ws.Range("A11").Value = 'CMG Champs Elysées - Boetie 7'
$ws.Range("B11").Value = 'US$14,434'
$ws.Range("C11").Value = '7.6'
$ws.Range("D11").Value = 'Good'
$ws.Range("E11").Value = '50'

$ws.Range("A12").Value = 'Rent a Room - Residence Meslay'
$ws.Range("B12").Value = 'US$5,750'
$ws.Range("C12").Value = '7.0'
$ws.Range("D12").Value = 'Good'
$ws.Range("E12").Value = '106'

$ws.Range("A13").Value = 'CMG Montorgueil X'
$ws.Range("B13").Value = 'US$6,753'
$ws.Range("C13").Value = '7.0'
$ws.Range("D13").Value = 'Good'
$ws.Range("E13").Value = '37'

$ws.Range("A14").Value = 'Rent a Room - Residence Blanche'
$ws.Range("B14").Value = 'US$8,456'
$ws.Range("C14").Value = '7.6'
$ws.Range("D14").Value = 'Good'
$ws.Range("E14").Value = '112'

$ws.Range("A15").Value = 'ARC de TRIOMPHE - FOCH PALACE'
$ws.Range("B15").Value = 'US$14,338'
$ws.Range("C15").Value = '8.8'
$ws.Range("D15").Value = 'Excellent'
$ws.Range("E15").Value = '18'

$ws.Range("A16").Value = 'CMG Jacques /Pantheon'
$ws.Range("B16").Value = 'US$11,275'
$ws.Range("C16").Value = '6.8'
$ws.Range("D16").Value = 'Review score'
$ws.Range("E16").Value = '32'

$ws.Range("A17").Value = 'CMG Résidence République II - Rue Béranger'
$ws.Range("B17").Value = 'US$15,296'
$ws.Range("C17").Value = '7.2'
$ws.Range("D17").Value = 'Good'
$ws.Range("E17").Value = '9'

$ws.Range("A18").Value = 'CMG Marais // République G'
$ws.Range("B18").Value = 'US$12,612'
$ws.Range("C18").Value = '8.1'
$ws.Range("D18").Value = 'Very Good'
$ws.Range("E18").Value = '59'

$ws.Range("A19").Value = 'Amazing apartment 8P3BDR - MontmartreSacré cœur'
$ws.Range("B19").Value = 'US$9,832'
$ws.Range("C19").Value = '8.6'
$ws.Range("D19").Value = 'Excellent'
$ws.Range("E19").Value = '21'

$ws.Range("A20").Value = 'Joyful apartment 2BR6P Heart of Paris - Louvre'
$ws.Range("B20").Value = 'US$11,560'
$ws.Range("C20").Value = '7.5'
$ws.Range("D20").Value = 'Good'
$ws.Range("E20").Value = '33'

$ws.Range("A21").Value = 'Luxury style appartement, Arc de Triomphe - Champs Elysées'
$ws.Range("B21").Value = 'US$20,334'
$ws.Range("C21").Value = '7.8'
$ws.Range("D21").Value = 'Good'
$ws.Range("E21").Value = '9'

$ws.Range("A22").Value = 'CMG Condorcet / Lentonnet I'
$ws.Range("B22").Value = 'US$16,335'
$ws.Range("C22").Value = '6.8'
$ws.Range("D22").Value = 'Review score'
$ws.Range("E22").Value = '10'

$ws.Range("A23").Value = 'Appartement Paris Tour Eiffel'
$ws.Range("B23").Value = 'US$24,446'
$ws.Range("C23").Value = '8.6'
$ws.Range("D23").Value = 'Excellent'
$ws.Range("E23").Value = '24'

$ws.Range("A24").Value = 'CMG - Nation Charonne G'
$ws.Range("B24").Value = 'US$11,860'
$ws.Range("C24").Value = '7.4'
$ws.Range("D24").Value = 'Good'
$ws.Range("E24").Value = '21'

$ws.Range("A25").Value = 'CMG - MAGENTA / SAINT-VINCENT DE PAUL'
$ws.Range("B25").Value = 'US$12,476'
$ws.Range("C25").Value = '7.1'
$ws.Range("D25").Value = 'Good'
$ws.Range("E25").Value = '16'

$ws.Range("A26").Value = 'CMG Place de la Nation I'
$ws.Range("B26").Value = 'US$13,057'
$ws.Range("C26").Value = '7.0'
$ws.Range("D26").Value = 'Good'
$ws.Range("E26").Value = '16'

$scoreRange.Style = "Normal"
$reviewsRange.Style = "Normal"

# Remove the two now-unused trailing rows (27 and 28).
$ws.Rows("27:28").Delete()
